# Update the cryptocurrency price/volume snapshot on Sheet1.
# Numeric-looking price strings are written with a leading apostrophe
# so Excel stores them as text (matching the original inlineStr cells)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.451.80'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '1.584.60'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''213.61'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '''44.38'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '''24.05'
$ws.Range("E10").Value = '  -2.08%  '
$ws.Range("E11").Value = '  -1.40%  '
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = '1.810.98'
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").Value = '1.589.81'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").Value = '''3.71'
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").Value = '''0.520'
$ws.Range("E16").Value = '  -1.55%  '
$ws.Range("D17").Value = '28.477.52'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = '''62.19'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").Value = '''230.56'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("E21").Value = '  -2.22%  '
$ws.Range("D23").Value = '''3.92'
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("E25").Value = '  +3.62%  '
$ws.Range("D26").Value = '''152.02'
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").Value = '''15.03'
$ws.Range("D28").Value = '''6.44'
$ws.Range("E28").Value = '  -1.49%  '
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").Value = '''0.0483'
$ws.Range("E31").Value = '  +2.75%  '
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("E34").Value = '  -2.27%  '
$ws.Range("D35").Value = '1.399.03'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("E36").Value = '  +6.27%  '
$ws.Range("E37").Value = '  -4.35%  '
$ws.Range("D38").Value = '''2.37'
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").Value = '''2.67'
$ws.Range("E39").Value = '  +1.97%  '
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''1.89'
$ws.Range("E43").Value = '  +1.64%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '''0.791'
$ws.Range("E44").Value = '  -2.46%  '
$ws.Range("D45").Value = '''0.0465'
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").Value = '''5.44'
$ws.Range("D47").Value = '''0.962'
$ws.Range("E47").Value = '  -1.93%  '
$ws.Range("D48").Value = '''63.17'
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("D49").Value = '1.722.23'
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = '''86.69'
$ws.Range("E51").Value = '  -2.32%  '